# Update workbook with the "http2 non ssl" re-run results.
# The "HTTP/2" and "HTTP/2 SSL" column headers are swapped on every sheet
# (column D is now "HTTP/2 SSL", column F is now "HTTP/2"), and the
# benchmark numbers for columns D, E, F, G (rows 4-7) are replaced with the
# newly measured values.

$wb = $excel.ActiveWorkbook

function Set-SheetValues {
    param(
        [string]$SheetName,
        [hashtable]$Values
    )

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($addr in $Values.Keys) {
        $ws.Range($addr).Value = $Values[$addr]
    }
}

Set-SheetValues "Transfer Time (s)" @{
    "D1" = "HTTP/2 SSL"
    "F1" = "HTTP/2"

    "D4" = 0.01860912831624349
    "E4" = 0.01781765832560506
    "F4" = 0.0108374330997467
    "G4" = 0.00646512990140088

    "D5" = 0.02502824942270915
    "E5" = 0.01429187471737145
    "F5" = 0.01645312070846558
    "G5" = 0.006651154588850086

    "D6" = 0.1433880964914958
    "E6" = 0.04356813924271643
    "F6" = 0.09403650760650635
    "G6" = 0.005134282398925803

    "D7" = 1.975765943527222
    "F7" = 0.7036097049713135
}

Set-SheetValues "Throughput (bps)" @{
    "D1" = "HTTP/2 SSL"
    "F1" = "HTTP/2"

    "D4" = 829991.1173650185
    "E4" = 394289.2181876033
    "F4" = 1187626.088236784
    "G4" = 529051.1059366131

    "D5" = 4802910.522485345
    "E5" = 1558273.501772276
    "F5" = 6910514.160241783
    "G5" = 2019480.877227425

    "D6" = 7891914.411074576
    "E6" = 2129848.297298706
    "F6" = 11181189.72082671
    "G6" = 622199.0804710645

    "D7" = 6392457.03444244
    "F7" = 14902807.51660114
}

Set-SheetValues "Overhead Ratio" @{
    "D1" = "HTTP/2 SSL"
    "F1" = "HTTP/2"

    "D4" = 1.021321614583333
    "E4" = 0
    "F4" = 1.02021611328125
    "G4" = 0.00004014610310760497

    "D5" = 1.00267578125
    "F5" = 1.002568359375

    "D6" = 1.000751495361328
    "F6" = 1.000741004943848

    "D7" = 1.000569629669189
    "F7" = 1.000568580627441
}

$wb.Save()
